$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet tracks one price "snapshot" per run in its own timestamped
# column, with two fixed columns ("nom" and "url_produit") always
# trailing at the end. This run adds a new snapshot column right
# before those two fixed columns.
#
# Concretely: insert a new column at FN. That shifts the old FN
# ("nom") to FO and the old FO ("url_produit") to FP, and grows the
# sheet dimension from A1:FO208 to A1:FP208 automatically.
$ws.Range("FN1").EntireColumn.Insert()

# Header for the newly inserted snapshot column: the timestamp of
# this run.
$ws.Cells.Item(1, 170).Value2 = "2026-02-04 19:32:04"

# The new snapshot simply repeats the previous snapshot (old FM,
# still column FM/169 after the insert) for every product row, since
# no new scrape happened between the two timestamps. Rows whose
# latest snapshot is blank are left blank too (the freshly inserted
# column already starts out empty, so there is nothing to do for
# those rows).
for ($r = 2; $r -le 208; $r++) {
    $latest = $ws.Cells.Item($r, 169).Value2
    if ($latest -ne $null -and $latest -ne "") {
        $ws.Cells.Item($r, 170).Value2 = $latest
    }
}
